# Applies the row-level corrections to the Fruta / Macroferia Regional de Talca - Granada dataset
# (Fruta / hortaliza, semanal refresh): dates, quality, volume, price and origin columns
# for rows 2-25 are re-populated with the corrected weekly values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44680
$ws.Cells.Item(2, 12).Value = 'Primera'
$ws.Cells.Item(2, 13).Value = 200
$ws.Cells.Item(2, 14).Value = 15000
$ws.Cells.Item(2, 15).Value = 15000
$ws.Cells.Item(2, 16).Value = 15000
$ws.Cells.Item(2, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(2, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(2, 19).Value = 1000
$ws.Cells.Item(2, 20).Value = 15

# Row 3
$ws.Cells.Item(3, 4).Value = 45085
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 280
$ws.Cells.Item(3, 14).Value = 17000
$ws.Cells.Item(3, 15).Value = 18000
$ws.Cells.Item(3, 16).Value = 17357
$ws.Cells.Item(3, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(3, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(3, 19).Value = 964
$ws.Cells.Item(3, 20).Value = 18

# Row 4
$ws.Cells.Item(4, 4).Value = 44355
$ws.Cells.Item(4, 12).Value = 'Especial'
$ws.Cells.Item(4, 13).Value = 50
$ws.Cells.Item(4, 14).Value = 18000
$ws.Cells.Item(4, 15).Value = 18000
$ws.Cells.Item(4, 16).Value = 18000
$ws.Cells.Item(4, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(4, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(4, 19).Value = 1000
$ws.Cells.Item(4, 20).Value = 18

# Row 5
$ws.Cells.Item(5, 4).Value = 44342
$ws.Cells.Item(5, 12).Value = 'Especial'
$ws.Cells.Item(5, 13).Value = 300
$ws.Cells.Item(5, 14).Value = 20000
$ws.Cells.Item(5, 15).Value = 20000
$ws.Cells.Item(5, 16).Value = 20000
$ws.Cells.Item(5, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(5, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(5, 19).Value = 1111
$ws.Cells.Item(5, 20).Value = 18

# Row 6
$ws.Cells.Item(6, 4).Value = 44714
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 100
$ws.Cells.Item(6, 14).Value = 20000
$ws.Cells.Item(6, 15).Value = 20000
$ws.Cells.Item(6, 16).Value = 20000
$ws.Cells.Item(6, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(6, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(6, 19).Value = 1111
$ws.Cells.Item(6, 20).Value = 18

# Row 7
$ws.Cells.Item(7, 4).Value = 45054
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 200
$ws.Cells.Item(7, 14).Value = 14000
$ws.Cells.Item(7, 15).Value = 14000
$ws.Cells.Item(7, 16).Value = 14000
$ws.Cells.Item(7, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(7, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(7, 19).Value = 778
$ws.Cells.Item(7, 20).Value = 18

# Row 8
$ws.Cells.Item(8, 4).Value = 45099
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 100
$ws.Cells.Item(8, 14).Value = 17000
$ws.Cells.Item(8, 15).Value = 17000
$ws.Cells.Item(8, 16).Value = 17000
$ws.Cells.Item(8, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(8, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(8, 19).Value = 944
$ws.Cells.Item(8, 20).Value = 18

# Row 9
$ws.Cells.Item(9, 4).Value = 44348
$ws.Cells.Item(9, 12).Value = 'Especial'
$ws.Cells.Item(9, 13).Value = 200
$ws.Cells.Item(9, 14).Value = 20000
$ws.Cells.Item(9, 15).Value = 20000
$ws.Cells.Item(9, 16).Value = 20000
$ws.Cells.Item(9, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(9, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(9, 19).Value = 1111
$ws.Cells.Item(9, 20).Value = 18

# Row 10
$ws.Cells.Item(10, 4).Value = 44319
$ws.Cells.Item(10, 12).Value = 'Especial'
$ws.Cells.Item(10, 13).Value = 120
$ws.Cells.Item(10, 14).Value = 20000
$ws.Cells.Item(10, 15).Value = 20000
$ws.Cells.Item(10, 16).Value = 20000
$ws.Cells.Item(10, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(10, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(10, 19).Value = 1111
$ws.Cells.Item(10, 20).Value = 18

# Row 11
$ws.Cells.Item(11, 4).Value = 44354
$ws.Cells.Item(11, 12).Value = 'Primera'
$ws.Cells.Item(11, 13).Value = 100
$ws.Cells.Item(11, 14).Value = 18000
$ws.Cells.Item(11, 15).Value = 18000
$ws.Cells.Item(11, 16).Value = 18000
$ws.Cells.Item(11, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(11, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(11, 19).Value = 1000
$ws.Cells.Item(11, 20).Value = 18

# Row 12
$ws.Cells.Item(12, 4).Value = 44326
$ws.Cells.Item(12, 12).Value = 'Especial'
$ws.Cells.Item(12, 13).Value = 300
$ws.Cells.Item(12, 14).Value = 20000
$ws.Cells.Item(12, 15).Value = 20000
$ws.Cells.Item(12, 16).Value = 20000
$ws.Cells.Item(12, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(12, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(12, 19).Value = 1111
$ws.Cells.Item(12, 20).Value = 18

# Row 13
$ws.Cells.Item(13, 4).Value = 44691
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 13).Value = 100
$ws.Cells.Item(13, 14).Value = 17000
$ws.Cells.Item(13, 15).Value = 17000
$ws.Cells.Item(13, 16).Value = 17000
$ws.Cells.Item(13, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(13, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(13, 19).Value = 944
$ws.Cells.Item(13, 20).Value = 18

# Row 14
$ws.Cells.Item(14, 4).Value = 44328
$ws.Cells.Item(14, 12).Value = 'Especial'
$ws.Cells.Item(14, 13).Value = 250
$ws.Cells.Item(14, 14).Value = 20000
$ws.Cells.Item(14, 15).Value = 20000
$ws.Cells.Item(14, 16).Value = 20000
$ws.Cells.Item(14, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(14, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(14, 19).Value = 1111
$ws.Cells.Item(14, 20).Value = 18

# Row 15
$ws.Cells.Item(15, 4).Value = 44294
$ws.Cells.Item(15, 12).Value = 'Primera'
$ws.Cells.Item(15, 13).Value = 50
$ws.Cells.Item(15, 14).Value = 12000
$ws.Cells.Item(15, 15).Value = 12000
$ws.Cells.Item(15, 16).Value = 12000
$ws.Cells.Item(15, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(15, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(15, 19).Value = 800
$ws.Cells.Item(15, 20).Value = 15

# Row 17
$ws.Cells.Item(17, 4).Value = 44299
$ws.Cells.Item(17, 12).Value = 'Primera'
$ws.Cells.Item(17, 13).Value = 100
$ws.Cells.Item(17, 14).Value = 15000
$ws.Cells.Item(17, 15).Value = 15000
$ws.Cells.Item(17, 16).Value = 15000
$ws.Cells.Item(17, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(17, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(17, 19).Value = 1000
$ws.Cells.Item(17, 20).Value = 15

# Row 19
$ws.Cells.Item(19, 4).Value = 45062
$ws.Cells.Item(19, 12).Value = 'Primera'
$ws.Cells.Item(19, 13).Value = 200
$ws.Cells.Item(19, 14).Value = 15000
$ws.Cells.Item(19, 15).Value = 15000
$ws.Cells.Item(19, 16).Value = 15000
$ws.Cells.Item(19, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(19, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(19, 19).Value = 1000
$ws.Cells.Item(19, 20).Value = 15

# Row 20
$ws.Cells.Item(20, 4).Value = 44291
$ws.Cells.Item(20, 12).Value = 'Primera'
$ws.Cells.Item(20, 13).Value = 150
$ws.Cells.Item(20, 14).Value = 12000
$ws.Cells.Item(20, 15).Value = 12000
$ws.Cells.Item(20, 16).Value = 12000
$ws.Cells.Item(20, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(20, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(20, 19).Value = 800
$ws.Cells.Item(20, 20).Value = 15

# Row 21
$ws.Cells.Item(21, 4).Value = 45083
$ws.Cells.Item(21, 12).Value = 'Primera'
$ws.Cells.Item(21, 13).Value = 120
$ws.Cells.Item(21, 14).Value = 17000
$ws.Cells.Item(21, 15).Value = 17000
$ws.Cells.Item(21, 16).Value = 17000
$ws.Cells.Item(21, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(21, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(21, 19).Value = 944
$ws.Cells.Item(21, 20).Value = 18

# Row 22
$ws.Cells.Item(22, 4).Value = 44340
$ws.Cells.Item(22, 12).Value = 'Primera'
$ws.Cells.Item(22, 13).Value = 230
$ws.Cells.Item(22, 14).Value = 20000
$ws.Cells.Item(22, 15).Value = 20000
$ws.Cells.Item(22, 16).Value = 20000
$ws.Cells.Item(22, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(22, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(22, 19).Value = 1111
$ws.Cells.Item(22, 20).Value = 18

# Row 23
$ws.Cells.Item(23, 4).Value = 44358
$ws.Cells.Item(23, 12).Value = 'Especial'
$ws.Cells.Item(23, 13).Value = 150
$ws.Cells.Item(23, 14).Value = 18000
$ws.Cells.Item(23, 15).Value = 18000
$ws.Cells.Item(23, 16).Value = 18000
$ws.Cells.Item(23, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(23, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(23, 19).Value = 1000
$ws.Cells.Item(23, 20).Value = 18

# Row 24
$ws.Cells.Item(24, 4).Value = 44358
$ws.Cells.Item(24, 12).Value = 'Primera'
$ws.Cells.Item(24, 13).Value = 100
$ws.Cells.Item(24, 14).Value = 17000
$ws.Cells.Item(24, 15).Value = 17000
$ws.Cells.Item(24, 16).Value = 17000
$ws.Cells.Item(24, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(24, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(24, 19).Value = 944
$ws.Cells.Item(24, 20).Value = 18

# Row 25
$ws.Cells.Item(25, 4).Value = 44692
$ws.Cells.Item(25, 12).Value = 'Especial'
$ws.Cells.Item(25, 13).Value = 150
$ws.Cells.Item(25, 14).Value = 17000
$ws.Cells.Item(25, 15).Value = 17000
$ws.Cells.Item(25, 16).Value = 17000
$ws.Cells.Item(25, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(25, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(25, 19).Value = 944
$ws.Cells.Item(25, 20).Value = 18
